$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-placeholder value for May hotel occupancy (C8)
$ws.Range("C8").Value = 15.3

# Move the active selection cursor (cosmetic, matches author's last position)
$ws.Range("F13").Select()
